# Update transition-probability matrix on Sheet1 with recalculated values
# (more games simulated / faster simulate-game logic / optimization logic draft).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2228915662650602
$ws.Range("C2").Value = 0.5060240963855421
$ws.Range("J2").Value = 0.02409638554216868
$ws.Range("P2").Value = 0.1656626506024096
$ws.Range("S2").Value = 0.08132530120481928
$ws.Range("B3").Value = 0.01796407185628742
$ws.Range("C3").Value = 0.005988023952095809
$ws.Range("J3").Value = 0.01197604790419162
$ws.Range("P3").Value = 0.7544910179640718
$ws.Range("S3").Value = 0.2095808383233533
$ws.Range("J4").Value = 0.1282051282051282
$ws.Range("P4").Value = 0.6923076923076923
$ws.Range("S4").Value = 0.1794871794871795
$ws.Range("B6").Value = 0.08695652173913043
$ws.Range("D6").Value = 0.01932367149758454
$ws.Range("E6").Value = 0.004830917874396135
$ws.Range("F6").Value = 0.07246376811594203
$ws.Range("J6").Value = 0.2318840579710145
$ws.Range("O6").Value = 0.00966183574879227
$ws.Range("Q6").Value = 0.1980676328502415
$ws.Range("R6").Value = 0.04347826086956522
$ws.Range("S6").Value = 0.3333333333333333
$ws.Range("B7").Value = 0.08247422680412371
$ws.Range("D7").Value = 0.02577319587628866
$ws.Range("F7").Value = 0.07731958762886598
$ws.Range("J7").Value = 0.1494845360824742
$ws.Range("O7").Value = 0.03092783505154639
$ws.Range("Q7").Value = 0.2010309278350516
$ws.Range("R7").Value = 0.05154639175257732
$ws.Range("S7").Value = 0.3814432989690721
$ws.Range("B8").Value = 0.1412742382271468
$ws.Range("D8").Value = 0.01662049861495845
$ws.Range("F8").Value = 0.05263157894736842
$ws.Range("J8").Value = 0.1163434903047091
$ws.Range("O8").Value = 0.02493074792243767
$ws.Range("Q8").Value = 0.2022160664819945
$ws.Range("R8").Value = 0.0443213296398892
$ws.Range("S8").Value = 0.4016620498614958
$ws.Range("B9").Value = 0.09236947791164658
$ws.Range("D9").Value = 0.008032128514056224
$ws.Range("E9").Value = 0.004016064257028112
$ws.Range("F9").Value = 0.05622489959839357
$ws.Range("J9").Value = 0.1004016064257028
$ws.Range("O9").Value = 0.02409638554216868
$ws.Range("Q9").Value = 0.1887550200803213
$ws.Range("R9").Value = 0.0963855421686747
$ws.Range("S9").Value = 0.429718875502008
$ws.Range("B10").Value = 0.1172413793103448
$ws.Range("D10").Value = 0.02068965517241379
$ws.Range("F10").Value = 0.06637931034482758
$ws.Range("J10").Value = 0.1405172413793103
$ws.Range("O10").Value = 0.02672413793103448
$ws.Range("Q10").Value = 0.2353448275862069
$ws.Range("R10").Value = 0.05517241379310345
$ws.Range("S10").Value = 0.3379310344827586
$ws.Range("G11").Value = 0.1324503311258278
$ws.Range("J11").Value = 0.08609271523178808
$ws.Range("K11").Value = 0.1655629139072848
$ws.Range("L11").Value = 0.5927152317880795
$ws.Range("S11").Value = 0.02317880794701987
$ws.Range("G12").Value = 0.7252747252747253
$ws.Range("J12").Value = 0.2307692307692308
$ws.Range("K12").Value = 0.005494505494505495
$ws.Range("L12").Value = 0.01098901098901099
$ws.Range("S12").Value = 0.02747252747252747
$ws.Range("G13").Value = 0.7291666666666666
$ws.Range("J13").Value = 0.2291666666666667
$ws.Range("S13").Value = 0.04166666666666666
$ws.Range("F15").Value = 0.03720930232558139
$ws.Range("H15").Value = 0.1674418604651163
$ws.Range("I15").Value = 0.08372093023255814
$ws.Range("J15").Value = 0.3534883720930233
$ws.Range("K15").Value = 0.07906976744186046
$ws.Range("M15").Value = 0.004651162790697674
$ws.Range("O15").Value = 0.04651162790697674
$ws.Range("S15").Value = 0.2279069767441861
$ws.Range("F16").Value = 0.03465346534653466
$ws.Range("H16").Value = 0.1831683168316832
$ws.Range("I16").Value = 0.0891089108910891
$ws.Range("J16").Value = 0.3762376237623762
$ws.Range("K16").Value = 0.1188118811881188
$ws.Range("O16").Value = 0.03465346534653466
$ws.Range("S16").Value = 0.1633663366336634
$ws.Range("F17").Value = 0.008474576271186441
$ws.Range("H17").Value = 0.1610169491525424
$ws.Range("I17").Value = 0.125
$ws.Range("J17").Value = 0.4088983050847458
$ws.Range("K17").Value = 0.08686440677966102
$ws.Range("M17").Value = 0.01483050847457627
$ws.Range("N17").Value = 0.00423728813559322
$ws.Range("O17").Value = 0.07627118644067797
$ws.Range("S17").Value = 0.1144067796610169
$ws.Range("F18").Value = 0.01626016260162602
$ws.Range("H18").Value = 0.1219512195121951
$ws.Range("I18").Value = 0.2032520325203252
$ws.Range("J18").Value = 0.3821138211382114
$ws.Range("K18").Value = 0.07317073170731707
$ws.Range("M18").Value = 0.04065040650406504
$ws.Range("O18").Value = 0.06504065040650407
$ws.Range("S18").Value = 0.0975609756097561
$ws.Range("F19").Value = 0.02378854625550661
$ws.Range("H19").Value = 0.1726872246696035
$ws.Range("I19").Value = 0.1145374449339207
$ws.Range("J19").Value = 0.345374449339207
$ws.Range("K19").Value = 0.1251101321585903
$ws.Range("M19").Value = 0.03259911894273128
$ws.Range("N19").Value = 0.000881057268722467
$ws.Range("O19").Value = 0.06607929515418502
$ws.Range("S19").Value = 0.118942731277533
